# Applies the cryptos.xlsx data refresh described in the commit:
# "Updated cryptos list on Thu May 16 05:54:15 UTC 2024 with GitHub Actions"
#
# The sheet stores Price (D) / Volume(1h) (E) as literal text (t="inlineStr")
# even when the text looks like a plain number (e.g. "6.72"). Excel's COM
# layer auto-coerces a bare numeric-looking string assigned via .Value into a
# real number, so for any new Price string that parses as a float we briefly
# force the cell to Text format, assign it, then restore the style (dropping
# the now-unreferenced explicit number format) so no stray style index sticks
# around on cells that never had one. Percent strings keep their padding
# spaces ("  +6.37%  "), which already blocks Excel's numeric auto-detect, so
# they're just assigned directly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Plain text / non-numeric-looking values: direct assignment ---
$ws.Range('B46').Value = 'Maker'
$ws.Range('B47').Value = 'VeChain'
$ws.Range('C46').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('C47').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D2').Value = '65.852.72'
$ws.Range('D3').Value = '3.009.99'
$ws.Range('D9').Value = '3.006.48'
$ws.Range('D16').Value = '65.829.54'
$ws.Range('D17').Value = '3.506.81'
$ws.Range('D19').Value = '3.009.36'
$ws.Range('D46').Value = '2.790.86'
$ws.Range('E2').Value = '  +6.37%  '
$ws.Range('E3').Value = '  +3.66%  '
$ws.Range('E4').Value = '  -0.09%  '
$ws.Range('E5').Value = '  +2.00%  '
$ws.Range('E6').Value = '  +12.72%  '
$ws.Range('E7').Value = '  -0.09%  '
$ws.Range('E8').Value = '  +3.63%  '
$ws.Range('E9').Value = '  +3.65%  '
$ws.Range('E10').Value = '  -4.06%  '
$ws.Range('E11').Value = '  +7.15%  '
$ws.Range('E12').Value = '  +7.49%  '
$ws.Range('E13').Value = '  +8.98%  '
$ws.Range('E14').Value = '  +7.91%  '
$ws.Range('E15').Value = '  -0.35%  '
$ws.Range('E17').Value = '  +3.60%  '
$ws.Range('E18').Value = '  +7.64%  '
$ws.Range('E19').Value = '  +3.59%  '
$ws.Range('E20').Value = '  +6.76%  '
$ws.Range('E21').Value = '  +8.50%  '
$ws.Range('E22').Value = '  +5.89%  '
$ws.Range('E23').Value = '  +7.54%  '
$ws.Range('E24').Value = '  +4.47%  '
$ws.Range('E25').Value = '  +12.90%  '
$ws.Range('E26').Value = '  +3.30%  '
$ws.Range('E27').Value = '  +6.02%  '
$ws.Range('E28').Value = '  -0.03%  '
$ws.Range('E29').Value = '  +17.13%  '
$ws.Range('E30').Value = '  +17.03%  '
$ws.Range('E31').Value = '  -5.39%  '
$ws.Range('E32').Value = '  +3.47%  '
$ws.Range('E33').Value = '  +5.59%  '
$ws.Range('E34').Value = '  +3.93%  '
$ws.Range('E35').Value = '  -0.11%  '
$ws.Range('E36').Value = '  +3.52%  '
$ws.Range('E37').Value = '  +8.35%  '
$ws.Range('E38').Value = '  +14.17%  '
$ws.Range('E39').Value = '  +2.02%  '
$ws.Range('E40').Value = '  +3.65%  '
$ws.Range('E41').Value = '  +16.40%  '
$ws.Range('E42').Value = '  +7.71%  '
$ws.Range('E43').Value = '  +6.81%  '
$ws.Range('E44').Value = '  +3.84%  '
$ws.Range('E45').Value = '  +13.18%  '
$ws.Range('E46').Value = '  +3.37%  '
$ws.Range('E47').Value = '  +5.85%  '
$ws.Range('E48').Value = '  +2.27%  '
$ws.Range('E50').Value = '  +10.66%  '
$ws.Range('E51').Value = '  +4.38%  '

# --- Numeric-looking Price strings: force Text format so they stay strings ---
$cell = $ws.Range('D5')
$cell.NumberFormat = '@'
$cell.Value = '581.66'
$cell.Style = 'Normal'
$cell = $ws.Range('D6')
$cell.NumberFormat = '@'
$cell.Value = '162.27'
$cell.Style = 'Normal'
$cell = $ws.Range('D8')
$cell.NumberFormat = '@'
$cell.Value = '0.519'
$cell.Style = 'Normal'
$cell = $ws.Range('D10')
$cell.NumberFormat = '@'
$cell.Value = '6.72'
$cell.Style = 'Normal'
$cell = $ws.Range('D12')
$cell.NumberFormat = '@'
$cell.Value = '0.460'
$cell.Style = 'Normal'
$cell = $ws.Range('D13')
$cell.NumberFormat = '@'
$cell.Value = '0.0000251'
$cell.Style = 'Normal'
$cell = $ws.Range('D14')
$cell.NumberFormat = '@'
$cell.Value = '34.74'
$cell.Style = 'Normal'
$cell = $ws.Range('D18')
$cell.NumberFormat = '@'
$cell.Value = '6.98'
$cell.Style = 'Normal'
$cell = $ws.Range('D20')
$cell.NumberFormat = '@'
$cell.Value = '457.85'
$cell.Style = 'Normal'
$cell = $ws.Range('D21')
$cell.NumberFormat = '@'
$cell.Value = '13.99'
$cell.Style = 'Normal'
$cell = $ws.Range('D22')
$cell.NumberFormat = '@'
$cell.Value = '0.689'
$cell.Style = 'Normal'
$cell = $ws.Range('D23')
$cell.NumberFormat = '@'
$cell.Value = '7.38'
$cell.Style = 'Normal'
$cell = $ws.Range('D24')
$cell.NumberFormat = '@'
$cell.Value = '82.39'
$cell.Style = 'Normal'
$cell = $ws.Range('D26')
$cell.NumberFormat = '@'
$cell.Value = '12.40'
$cell.Style = 'Normal'
$cell = $ws.Range('D27')
$cell.NumberFormat = '@'
$cell.Value = '10.76'
$cell.Style = 'Normal'
$cell = $ws.Range('D29')
$cell.NumberFormat = '@'
$cell.Value = '8.15'
$cell.Style = 'Normal'
$cell = $ws.Range('D30')
$cell.NumberFormat = '@'
$cell.Value = '2.36'
$cell.Style = 'Normal'
$cell = $ws.Range('D33')
$cell.NumberFormat = '@'
$cell.Value = '26.99'
$cell.Style = 'Normal'
$cell = $ws.Range('D36')
$cell.NumberFormat = '@'
$cell.Value = '0.991'
$cell.Style = 'Normal'
$cell = $ws.Range('D38')
$cell.NumberFormat = '@'
$cell.Value = '2.17'
$cell.Style = 'Normal'
$cell = $ws.Range('D39')
$cell.NumberFormat = '@'
$cell.Value = '49.73'
$cell.Style = 'Normal'
$cell = $ws.Range('D40')
$cell.NumberFormat = '@'
$cell.Value = '2.98'
$cell.Style = 'Normal'
$cell = $ws.Range('D42')
$cell.NumberFormat = '@'
$cell.Value = '44.06'
$cell.Style = 'Normal'
$cell = $ws.Range('D44')
$cell.NumberFormat = '@'
$cell.Value = '8.46'
$cell.Style = 'Normal'
$cell = $ws.Range('D45')
$cell.NumberFormat = '@'
$cell.Value = '392.24'
$cell.Style = 'Normal'
$cell = $ws.Range('D47')
$cell.NumberFormat = '@'
$cell.Value = '0.0355'
$cell.Style = 'Normal'
$cell = $ws.Range('D48')
$cell.NumberFormat = '@'
$cell.Value = '134.74'
$cell.Style = 'Normal'
$cell = $ws.Range('D50')
$cell.NumberFormat = '@'
$cell.Value = '23.83'
$cell.Style = 'Normal'
